# Cigna_Global_Health/userType.xlsx
#
# Insert a new benefit row "Repatriation-Benefit" / "All" directly below the
# existing "Emergency-Evacuation" / "All" row (row 34), pushing every row
# below it down by one (old rows 35-39 become 36-40). The inserted row
# naturally inherits the formatting (cell style + row height) of the row
# above it, exactly like Excel does when a row is inserted through the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (shifts old rows 35-39 down to 36-40).
$ws.Rows(35).Insert() | Out-Null

# Populate the new row's two cells.
$ws.Range("A35").Value = "Repatriation-Benefit"
$ws.Range("B35").Value = "All"

# Reflect the cursor/selection state recorded in the edit.
$ws.Range("B34").Select() | Out-Null
